$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell values -----------------------------------------------------
$ws.Range("A1").Value = "Farmer Table"

$ws.Range("A2").Value = "Property Name"
$ws.Range("B2").Value = "Data Type"
$ws.Range("C2").Value = "Key"
$ws.Range("D2").Value = "Default"

$ws.Range("A3").Value = "Farmer Id"
$ws.Range("B3").Value = "Int"
$ws.Range("C3").Value = "Primary Key"
$ws.Range("D3").Value = "Not Null"

$ws.Range("A4").Value = "First Name"
$ws.Range("B4").Value = "Varchar"

$ws.Range("A5").Value = "Last Name"
$ws.Range("B5").Value = "Varchar"

$ws.Range("A6").Value = "Phone"
$ws.Range("B6").Value = "Big Int"

$ws.Range("A7").Value = "Mail Id"
$ws.Range("B7").Value = "Varchar"

$ws.Range("A8").Value = "Address"
$ws.Range("B8").Value = "Varchar"

# --- Column widths -----------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 22
$ws.Columns.Item(2).ColumnWidth = 19.6328125
$ws.Columns.Item(3).ColumnWidth = 17.81640625
$ws.Columns.Item(4).ColumnWidth = 16.1796875

# --- Merge the title row first (before any border exists) so the merge
#     does not split the surrounding border into per-edge variants -----
$ws.Range("A1:D1").Merge()

# --- Borders around the whole table ------------------------------------
$ws.Range("A1:D8").Borders.LineStyle = 1

# --- Bold header rows ----------------------------------------------------
$ws.Range("A2:D2").Font.Bold = $true
$ws.Range("A1:D1").Font.Bold = $true
$ws.Range("A1:D1").HorizontalAlignment = -4108

# --- Page setup ------------------------------------------------------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- Selection ------------------------------------------------------
$ws.Range("B9").Select()
